$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Remove the 3rd example row for Body.search (row 8), shifting all rows below up by one.
$ws.Rows.Item(8).Delete()

# Restore selection/view state similar to the authored change.
$ws.Activate()
$ws.Range("B7").Select()
